# Insert a new data row above the current row 323 ("Feria Lagunitas de
# Puerto Montt" / Apio sheet). This shifts the existing rows 323..411 down
# to 324..412 (and the sheet dimension grows from R411 to R412), then fills
# the newly-inserted row 323 with its own record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("323:323").Insert()

$ws.Cells.Item(323, 1).Value  = 4
$ws.Cells.Item(323, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(323, 3).Value  = 'Los Lagos'
$ws.Cells.Item(323, 4).Value  = 44985
$ws.Cells.Item(323, 5).Value  = 10
$ws.Cells.Item(323, 6).Value  = 100112017
$ws.Cells.Item(323, 7).Value  = 'Apio'
$ws.Cells.Item(323, 8).Value  = 'Americana (o)'
$ws.Cells.Item(323, 9).Value  = 'Segunda'
$ws.Cells.Item(323, 10).Value = 50
$ws.Cells.Item(323, 11).Value = 10000
$ws.Cells.Item(323, 12).Value = 11000
$ws.Cells.Item(323, 13).Value = 10500
$ws.Cells.Item(323, 14).Value = '$/docena de matas'
$ws.Cells.Item(323, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(323, 16).Value = 1750
$ws.Cells.Item(323, 17).Value = 6
$ws.Cells.Item(323, 18).Value = 'Hortaliza'
